$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Median Value" (column C) scores, recalculated relative to the median
# AFTER merging with zip/census tract data, and the corresponding "Tier"
# (column D) labels that shift as a result.
$updates = @{
    2 = @{ C = 1.028985507246377; D = "4th Tier" }
    3 = @{ C = 0.5090579710144928; D = $null }
    4 = @{ C = 1.957971014492754; D = $null }
    5 = @{ C = 1.318840579710145; D = "3rd Tier" }
    6 = @{ C = 1.540760869565217; D = $null }
    7 = @{ C = 1; D = $null }
    8 = @{ C = 0.6898550724637681; D = $null }
    9 = @{ C = 1.855676328502415; D = $null }
    10 = @{ C = 0.7608695652173912; D = $null }
    11 = @{ C = 1.00054347826087; D = $null }
    12 = @{ C = 1.032608695652174; D = "4th Tier" }
    13 = @{ C = 1.087409420289855; D = "4th Tier" }
    14 = @{ C = 1.334692028985507; D = $null }
    15 = @{ C = 1.458333333333333; D = "2nd Tier" }
    16 = @{ C = 0.8293478260869566; D = $null }
    17 = @{ C = 0.9146286231884058; D = "Below Median" }
    18 = @{ C = 0.5727657004830917; D = $null }
    19 = @{ C = 0.8510466988727858; D = "Below Median" }
    20 = @{ C = 0.7059178743961352; D = $null }
    21 = @{ C = 0.6518115942028985; D = $null }
    22 = @{ C = 0.5217391304347826; D = $null }
    23 = @{ C = 0.6105072463768116; D = $null }
    24 = @{ C = 1.389855072463768; D = $null }
    25 = @{ C = 1.675724637681159; D = $null }
    26 = @{ C = 1.361111111111111; D = $null }
    27 = @{ C = 1.499547101449275; D = $null }
    28 = @{ C = 1.001811594202898; D = $null }
    29 = @{ C = 0.5757246376811593; D = $null }
    30 = @{ C = 1.43677536231884; D = $null }
    31 = @{ C = 1.059581320450886; D = "4th Tier" }
    32 = @{ C = 1.678985507246377; D = $null }
    33 = @{ C = 1.016606280193237; D = "4th Tier" }
    34 = @{ C = 0.9710144927536231; D = "Below Median" }
    35 = @{ C = 0.4839975845410627; D = $null }
    36 = @{ C = 1.356884057971014; D = $null }
    37 = @{ C = 0.7355072463768115; D = $null }
    38 = @{ C = 1.27536231884058; D = "3rd Tier" }
    39 = @{ C = 1.306159420289855; D = "3rd Tier" }
    40 = @{ C = 1.151449275362319; D = $null }
    41 = @{ C = 1.109601449275362; D = "4th Tier" }
    42 = @{ C = 0.8327294685990339; D = $null }
    43 = @{ C = 1.534420289855072; D = $null }
    44 = @{ C = 0.8876811594202898; D = "Below Median" }
    45 = @{ C = 1.113405797101449; D = $null }
    46 = @{ C = 0.601086956521739; D = $null }
    47 = @{ C = 0.9658816425120773; D = "Below Median" }
    48 = @{ C = 1.185688405797101; D = $null }
    49 = @{ C = 1.41268115942029; D = $null }
    50 = @{ C = 1.071557971014493; D = "4th Tier" }
    51 = @{ C = 0.8834541062801933; D = "Below Median" }
    52 = @{ C = 0.6644927536231884; D = $null }
    53 = @{ C = 1.204710144927536; D = $null }
    54 = @{ C = 0.9035326086956521; D = "Below Median" }
    55 = @{ C = 1.123188405797101; D = $null }
    56 = @{ C = 0.9184782608695652; D = "Below Median" }
    57 = @{ C = 0.527536231884058; D = $null }
    58 = @{ C = 0.6124999999999999; D = $null }
    59 = @{ C = 0.3043478260869565; D = $null }
    60 = @{ C = 0.5126811594202898; D = $null }
    61 = @{ C = 0.8211050724637681; D = $null }
    62 = @{ C = 1.341032608695652; D = $null }
    63 = @{ C = 0.6391304347826087; D = $null }
    64 = @{ C = 0.5807971014492753; D = $null }
    65 = @{ C = 0.6625905797101449; D = $null }
    66 = @{ C = 0.3719806763285024; D = $null }
    67 = @{ C = 0.9739130434782608; D = "Below Median" }
    68 = @{ C = 1.911684782608696; D = $null }
    69 = @{ C = 0.6036231884057971; D = $null }
    70 = @{ C = 0.8876811594202898; D = "Below Median" }
    71 = @{ C = 1.607971014492753; D = $null }
    72 = @{ C = 1.557246376811594; D = $null }
    73 = @{ C = 0.5427536231884057; D = $null }
    74 = @{ C = 1.327898550724637; D = "3rd Tier" }
    75 = @{ C = 1.296014492753623; D = "3rd Tier" }
    76 = @{ C = 1.348429951690821; D = $null }
    77 = @{ C = 0.7684782608695652; D = $null }
    78 = @{ C = 1.480525362318841; D = "2nd Tier" }
    79 = @{ C = 0.6657608695652174; D = $null }
    80 = @{ C = 0.6467391304347826; D = $null }
}

foreach ($row in $updates.Keys) {
    $entry = $updates[$row]
    $ws.Cells.Item($row, 3).Value = $entry.C
    if ($null -ne $entry.D) {
        $ws.Cells.Item($row, 4).Value = $entry.D
    }
}
